{"js": "const body = context.document.body;\n\n// Find the (first, in document order) paragraph whose entire text is the\n// literal placeholder \"(NEEDS RESPONSE)\".\nasync function getNeedsResponseParagraph() {\n  const results = body.search(\"(NEEDS RESPONSE)\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const hit = results.items[0];\n  const para = hit.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n  return para;\n}\n\n// --- Hunk 1: the \"(NEEDS RESPONSE)\" paragraph right after \"The results\n//     section is however disconcerting...\" becomes the reviewers' response,\n//     switched from italic to bold. ---\nconst para1 = await getNeedsResponseParagraph();\npara1.clear();\nawait context.sync();\nconst run1 = para1.insertText(\n  \"We believe that the responses to specific concerns (detailed below) help to clarify the message. This section is intended as a straightforward, objective assessment of the data available through ForC and submitted to EFDB, and is not intended to make a particular point.\",\n  Word.InsertLocation.replace\n);\nrun1.font.bold = true;\nrun1.font.boldBidirectional = true;\nawait context.sync();\n\n// --- Hunk 2: the \"(NEEDS RESPONSE)\" paragraph right after \"...Overall the\n//     discussion seems somewhat verbose.\" becomes a bold response, and a\n//     brand-new italic paragraph \"(VERBOSITY CRITICISM NEEDS RESPONSE)\" is\n//     inserted immediately after it. ---\nconst para2 = await getNeedsResponseParagraph();\n\n// Insert the new paragraph *before* touching para2's own text/formatting,\n// so the freshly-inserted paragraph inherits para2's still-italic run\n// formatting, exactly like the original text had.\nconst newPara = para2.insertParagraph(\"(VERBOSITY CRITICISM NEEDS RESPONSE)\", Word.InsertLocation.after);\nnewPara.font.italic = true;\nnewPara.font.italicBidirectional = true;\nawait context.sync();\n\npara2.clear();\nawait context.sync();\nconst run2 = para2.insertText(\n  \"We believe that the changes to the results section and figures now convincingly demonstrate the claims in the discussion.\",\n  Word.InsertLocation.replace\n);\nrun2.font.bold = true;\nrun2.font.boldBidirectional = true;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Helper: locate the next paragraph (at or after 1-based index $afterIndex)\n# whose text is exactly \"(NEEDS RESPONSE)\" and return a Range covering just\n# the visible text (the trailing paragraph-mark character is excluded).\nfunction Get-NeedsResponseRange($afterIndex) {\n    $count = $d.Paragraphs.Count\n    for ($i = $afterIndex; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $full = $p.Range\n        $trimmed = $d.Range($full.Start, $full.End - 1)\n        if ($trimmed.Text -eq \"(NEEDS RESPONSE)\") {\n            return $trimmed\n        }\n    }\n    return $null\n}\n\n# Helper: the Paragraph object whose Range starts right after $pos\n# (i.e. the paragraph immediately following the one ending at $pos).\nfunction Get-ParagraphStartingAt($pos) {\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Start -eq $pos) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# --- Hunk 1: first \"(NEEDS RESPONSE)\" (after \"The results section is\n#     however disconcerting...\") becomes the reviewer-response paragraph,\n#     switched from italic to bold. ---\n$range1 = Get-NeedsResponseRange 1\n$range1.Delete()\n$range1.InsertAfter(\"We believe that the responses to specific concerns (detailed below) help to clarify the message. This section is intended as a straightforward, objective assessment of the data available through ForC and submitted to EFDB, and is not intended to make a particular point.\")\n$range1.Bold = 1\n$range1.BoldBi = 1\n\n# --- Hunk 2: second \"(NEEDS RESPONSE)\" (after \"...Overall the discussion\n#     seems somewhat verbose.\") becomes a bold response paragraph, and a\n#     new italic paragraph \"(VERBOSITY CRITICISM NEEDS RESPONSE)\" is\n#     inserted right after it. ---\n$range2 = Get-NeedsResponseRange 1\n\n# Insert the new paragraph mark *before* changing $range2's text/format,\n# so the fresh (still-empty) paragraph inherits the original italic\n# run formatting, exactly like the source run had.\n$range2.InsertParagraphAfter()\n\n$range2.Delete()\n$range2.InsertAfter(\"We believe that the changes to the results section and figures now convincingly demonstrate the claims in the discussion.\")\n$range2.Bold = 1\n$range2.BoldBi = 1\n\n# The newly-inserted (still italic, still empty) paragraph immediately\n# follows the paragraph we just edited.\n$newPara = Get-ParagraphStartingAt ($range2.End + 1)\n$newFull = $newPara.Range\n$newTrimmed = $d.Range($newFull.Start, $newFull.End - 1)\n$newTrimmed.InsertAfter(\"(VERBOSITY CRITICISM NEEDS RESPONSE)\")\n"}
